# Correct misleading values in mapping schemes: revised area and cost
# assumptions for all occupancies; revised count assumptions for
# non-residential, on the "Dwellings_buildings" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# --- H:1 block (rows 2-4) ---
$ws.Range("C2").Value = 130
$ws.Range("D2").Value = 417.8037459283388

$ws.Range("B3").Value = "Professional and technical services"
$ws.Range("C3").Value = 130
$ws.Range("D3").Value = 479.899348534202

$ws.Range("B4").Value = "All other services"
$ws.Range("C4").Value = 130
$ws.Range("D4").Value = 433.1359934853421

# --- H:2 block (rows 5-7) ---
$ws.Range("C5").Value = 260
$ws.Range("D5").Value = 417.8037459283388
$ws.Range("F5").Value = 2

$ws.Range("B6").Value = "Professional and technical services"
$ws.Range("C6").Value = 260
$ws.Range("D6").Value = 479.899348534202
$ws.Range("F6").Value = 2

$ws.Range("B7").Value = "All other services"
$ws.Range("C7").Value = 260
$ws.Range("D7").Value = 433.1359934853421
$ws.Range("F7").Value = 2

# --- H:3 block (rows 8-10) ---
$ws.Range("C8").Value = 450
$ws.Range("D8").Value = 417.8037459283388
$ws.Range("F8").Value = 3

$ws.Range("B9").Value = "Professional and technical services"
$ws.Range("C9").Value = 450
$ws.Range("D9").Value = 479.899348534202
$ws.Range("F9").Value = 3

$ws.Range("B10").Value = "All other services"
$ws.Range("C10").Value = 450
$ws.Range("D10").Value = 433.1359934853421
$ws.Range("F10").Value = 3

# --- HBET:3-6 block (rows 11-13) ---
$ws.Range("C11").Value = 900
$ws.Range("D11").Value = 417.8037459283388
$ws.Range("F11").Value = 5

$ws.Range("B12").Value = "Professional and technical services"
$ws.Range("C12").Value = 900
$ws.Range("D12").Value = 479.899348534202
$ws.Range("F12").Value = 5

$ws.Range("B13").Value = "All other services"
$ws.Range("C13").Value = 900
$ws.Range("D13").Value = 433.1359934853421
$ws.Range("F13").Value = 5

# --- HBET:4-7 block (rows 14-16) ---
$ws.Range("B14").Value = "Professional and technical services"
$ws.Range("C14").Value = 1200
$ws.Range("D14").Value = 628.6221498371335
$ws.Range("F14").Value = 5

$ws.Range("B15").Value = "Professional and technical services"
$ws.Range("C15").Value = 1200
$ws.Range("D15").Value = 417.8037459283388
$ws.Range("F15").Value = 5

$ws.Range("B16").Value = "All other services"
$ws.Range("C16").Value = 1200
$ws.Range("D16").Value = 433.1359934853421
$ws.Range("F16").Value = 5

# --- HBET:8+ block (rows 17-19) ---
$ws.Range("B17").Value = "Professional and technical services"
$ws.Range("C17").Value = 3200
$ws.Range("D17").Value = 628.6221498371335
$ws.Range("F17").Value = 10

$ws.Range("B18").Value = "Professional and technical services"
$ws.Range("C18").Value = 3200
$ws.Range("D18").Value = 417.8037459283388
$ws.Range("F18").Value = 10

$ws.Range("B19").Value = "All other services"
$ws.Range("C19").Value = 3200
$ws.Range("D19").Value = 433.1359934853421
$ws.Range("F19").Value = 10
